# Refresh the cryptocurrency Price (D) and Volume(1h) (E) columns
# with the latest values from the GitHub Actions scrape.
#
# D-column values that parse as plain numbers are written with a
# leading apostrophe so Excel stores them as text (preserving
# thousands-dot formatting / trailing zeros / leading zeros) instead
# of silently coercing them to a Double.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.859.16"
$ws.Range("E2").Value = "  -0.19%  "
$ws.Range("D3").Value = "3.097.24"
$ws.Range("E3").Value = "  -0.31%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").Value = "'576.14"
$ws.Range("E5").Value = "  -0.34%  "
$ws.Range("D6").Value = "'177.09"
$ws.Range("E6").Value = "  +2.47%  "
$ws.Range("E7").Value = "  +0.09%  "
$ws.Range("D8").Value = "3.095.17"
$ws.Range("E8").Value = "  -0.26%  "
$ws.Range("D9").Value = "'0.509"
$ws.Range("E9").Value = "  -2.03%  "
$ws.Range("D10").Value = "'6.33"
$ws.Range("E10").Value = "  -2.26%  "
$ws.Range("E11").Value = "  -1.85%  "
$ws.Range("E12").Value = "  -2.77%  "
$ws.Range("D13").Value = "'0.0000240"
$ws.Range("E13").Value = "  -3.12%  "
$ws.Range("D14").Value = "'35.85"
$ws.Range("E14").Value = "  -2.53%  "
$ws.Range("D16").Value = "3.617.31"
$ws.Range("E16").Value = "  -0.13%  "
$ws.Range("D17").Value = "66.883.06"
$ws.Range("E17").Value = "  -0.11%  "
$ws.Range("D18").Value = "'6.95"
$ws.Range("E18").Value = "  -1.86%  "
$ws.Range("D19").Value = "'16.86"
$ws.Range("E19").Value = "  +2.64%  "
$ws.Range("D20").Value = "3.104.62"
$ws.Range("E20").Value = "  -0.13%  "
$ws.Range("D21").Value = "'481.30"
$ws.Range("E21").Value = "  -1.98%  "
$ws.Range("D22").Value = "'7.70"
$ws.Range("E22").Value = "  -2.20%  "
$ws.Range("D23").Value = "'0.685"
$ws.Range("E23").Value = "  -2.54%  "
$ws.Range("D24").Value = "'83.42"
$ws.Range("E24").Value = "  -0.50%  "
$ws.Range("D25").Value = "'12.58"
$ws.Range("E25").Value = "  -3.88%  "
$ws.Range("D26").Value = "'2.21"
$ws.Range("E26").Value = "  -2.98%  "
$ws.Range("D27").Value = "'10.17"
$ws.Range("E27").Value = "  -2.68%  "
$ws.Range("E28").Value = "  -0.01%  "
$ws.Range("D29").Value = "'8.01"
$ws.Range("E29").Value = "  +1.41%  "
$ws.Range("D30").Value = "'2.27"
$ws.Range("E30").Value = "  -4.27%  "
$ws.Range("E31").Value = "  -2.74%  "
$ws.Range("D32").Value = "'27.96"
$ws.Range("E32").Value = "  -1.20%  "
$ws.Range("E33").Value = "  -2.14%  "
$ws.Range("D34").Value = "0.0₃0938"
$ws.Range("E34").Value = "  -0.64%  "
$ws.Range("E35").Value = "  +0.11%  "
$ws.Range("D36").Value = "'48.40"
$ws.Range("E36").Value = "  +2.54%  "
$ws.Range("D37").Value = "'5.57"
$ws.Range("E37").Value = "  -4.50%  "
$ws.Range("D38").Value = "'0.939"
$ws.Range("E38").Value = "  -3.30%  "
$ws.Range("D39").Value = "'49.00"
$ws.Range("E39").Value = "  -2.03%  "
$ws.Range("D40").Value = "'0.308"
$ws.Range("E40").Value = "  -0.16%  "
$ws.Range("D41").Value = "'0.122"
$ws.Range("E41").Value = "  -0.37%  "
$ws.Range("D42").Value = "'1.97"
$ws.Range("E42").Value = "  -3.18%  "
$ws.Range("D43").Value = "'8.24"
$ws.Range("E43").Value = "  -2.39%  "
$ws.Range("D44").Value = "'2.66"
$ws.Range("E44").Value = "  +2.98%  "
$ws.Range("D45").Value = "2.785.83"
$ws.Range("E45").Value = "  -0.69%  "
$ws.Range("D46").Value = "'370.86"
$ws.Range("E46").Value = "  -3.44%  "
$ws.Range("D47").Value = "'135.41"
$ws.Range("E47").Value = "  -0.05%  "
$ws.Range("D48").Value = "'0.0342"
$ws.Range("E48").Value = "  -2.25%  "
$ws.Range("D50").Value = "'24.85"
$ws.Range("E50").Value = "  +0.33%  "
$ws.Range("D51").Value = "'2.21"
$ws.Range("E51").Value = "  +1.18%  "
